$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price & 1h volume change) per latest scrape.
# Numeric-looking price strings are apostrophe-prefixed to keep them as text
# (matching the source inlineStr cells), then the cell style is reset to
# "Normal" so no stray quote-prefix number-format style is left behind.

$ws.Range("D2").Value = "59.720.92"
$ws.Range("E2").Value = "  +8.28%  "
$ws.Range("D3").Value = "2.581.70"
$ws.Range("E3").Value = "  +10.14%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'506.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.57%  "
$ws.Range("D6").Value = "'157.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.90%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'0.611"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.35%  "
$ws.Range("D9").Value = "2.577.86"
$ws.Range("E9").Value = "  +9.83%  "
$ws.Range("D10").Value = "'6.11"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +12.50%  "
$ws.Range("E11").Value = "  +7.08%  "
$ws.Range("D12").Value = "'0.342"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.25%  "
$ws.Range("E13").Value = "  +1.36%  "
$ws.Range("D14").Value = "3.016.80"
$ws.Range("E14").Value = "  +9.65%  "
$ws.Range("D15").Value = "59.510.24"
$ws.Range("E15").Value = "  +7.98%  "
$ws.Range("D16").Value = "'21.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +9.58%  "
$ws.Range("E17").Value = "  +6.35%  "
$ws.Range("D18").Value = "2.577.27"
$ws.Range("E18").Value = "  +10.05%  "
$ws.Range("D19").Value = "'4.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.86%  "
$ws.Range("D20").Value = "'338.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.74%  "
$ws.Range("D21").Value = "'10.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +8.35%  "
$ws.Range("D22").Value = "'6.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.86%  "
$ws.Range("E23").Value = "  +0.44%  "
$ws.Range("D24").Value = "'60.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.46%  "
$ws.Range("D25").Value = "'0.418"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.12%  "
$ws.Range("E26").Value = "  +9.05%  "
$ws.Range("D27").Value = "2.679.66"
$ws.Range("E27").Value = "  +9.65%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").Value = "0.0₃0833"
$ws.Range("E29").Value = "  +12.20%  "
$ws.Range("E30").Value = "  +4.46%  "
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("D32").Value = "'156.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.02%  "
$ws.Range("D33").Value = "'19.43"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.06%  "
$ws.Range("E34").Value = "  +6.93%  "
$ws.Range("D35").Value = "'5.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.58%  "
$ws.Range("E36").Value = "  +9.61%  "
$ws.Range("D37").Value = "'3.91"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.35%  "
$ws.Range("D38").Value = "'0.859"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.28%  "
$ws.Range("D39").Value = "'303.18"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +21.11%  "
$ws.Range("D40").Value = "'3.74"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.90%  "
$ws.Range("E42").Value = "  +4.54%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").Value = "'0.102"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("E44").Value = "  +9.44%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "'0.0572"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.87%  "
$ws.Range("E46").Value = "  +25.66%  "
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("D48").Value = "'4.95"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +14.60%  "
$ws.Range("D49").Value = "'19.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +16.15%  "
$ws.Range("E50").Value = "  +7.65%  "
$ws.Range("E51").Value = "  +1.03%  "
